$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.585.26'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '2.494.71'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.31'
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.23'
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.512'
$ws.Range('E8').Value = '  -1.72%  '
$ws.Range('D9').Value = '2.493.03'
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.159'
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.355'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = '2.946.21'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').Value = '69.452.65'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.32'
$ws.Range('E17').Value = '  -2.73%  '
$ws.Range('D18').Value = '2.489.86'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.18'
$ws.Range('E19').Value = '  -2.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.39'
$ws.Range('E20').Value = '  -5.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '346.41'
$ws.Range('E21').Value = '  -1.29%  '
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('E23').Value = '  -3.19%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '69.68'
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('E26').Value = '  -2.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.68'
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('D28').Value = '2.616.42'
$ws.Range('E28').Value = '  -1.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').Value = '0.0₃0875'
$ws.Range('E30').Value = '  -2.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.66'
$ws.Range('E31').Value = '  -2.82%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '439.56'
$ws.Range('E32').Value = '  -5.44%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.19'
$ws.Range('E33').Value = '  -5.37%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  -2.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '155.21'
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('E37').Value = '  -4.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.07'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.18'
$ws.Range('E39').Value = '  -2.45%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.315'
$ws.Range('E41').Value = '  -1.33%  '
$ws.Range('E42').Value = '  -3.70%  '
$ws.Range('E43').Value = '  -2.04%  '
$ws.Range('E44').Value = '  -5.28%  '
$ws.Range('E45').Value = '  -6.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '138.74'
$ws.Range('E46').Value = '  -2.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.44'
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.512'
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0724'
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.573'
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('B51').Value = 'POPCAT'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sLBuDEsp6+popcat-popcat'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.95'
$ws.Range('E51').Value = '  +20.57%  '
